$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @{
    9  = -22.08130000000001
    13 = -22.2899
    16 = -21.7617
    18 = -22.0598
    20 = -21.72019999999997
    26 = -21.10699999999997
    27 = -21.9285
    29 = -20.97759999999997
    35 = -21.64159999999998
    36 = -20.80599999999999
    45 = -21.44769999999999
    55 = -22.1644
    57 = -22.379
    69 = -21.66639999999999
    76 = -19.61549999999999
    78 = -19.80069999999998
    82 = -21.646
    83 = -21.62329999999999
    93 = -21.41860000000002
    97 = -21.53679999999998
}

foreach ($row in $updates.Keys) {
    $ws.Range("A$row").Value = $updates[$row]
}
